# Generate Report for Handoff
# The localization run completed translation and the report is regenerated
# to reflect the new "Ready for handoff" status plus refreshed timestamps
# for the handoff/generation datetimes. Updating the Status text makes the
# Status column wider, so we also resize it (as AutoFit would).

$wb = $excel.ActiveWorkbook

# --- "Overview" sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-07 10:56:00"

# Status text got longer ("In Translation" -> "Ready for handoff"), widen
# the zh-cn / de-de status columns to fit.
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# --- "zh-cn" sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-07 10:55:55"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

# --- "de-de" sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-07 10:56:00"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
